$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37").Value = "Savignagne ✨"
$ws.Range("B37").Value = "ELIA BATTISTI | U.S. Guarna"
$ws.Range("C37").Value = "Andrea Conzatti | FC Savignano"
$ws.Range("D37").Value = "Stefano  Galvagni | Clitoriders"
$ws.Range("E37").Value = "Matteo Mazzola | MediaserT"
$ws.Range("F37").Value = "Simone Miorelli | SBARX"
